$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual "missing data" cell values (swap which cells are blank) ---
$ws.Range("C3").Value = 11.2
$ws.Range("E4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E9").Value = -6.8
$ws.Range("E10").Value = -6.1
$ws.Range("E11").Value = -7.9
$ws.Range("E12").Value = -5.3
$ws.Range("E15").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("C21").Value = 12.7
$ws.Range("C23").ClearContents()
$ws.Range("E33").Value = -8.1
$ws.Range("C34").Value = 10.5
$ws.Range("E34").Value = -6.4

# --- Remove the two rows that were dropped from the data set ---
# Delete the higher-numbered row first so the lower row index is unaffected.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

Write-Output "done"
